$wb = $excel.ActiveWorkbook

# --- 1. Extend the "Cuentas" table with a new "Tipo" column, all rows = "balance" ---
$wsCuentas = $wb.Worksheets.Item("Cuentas")
$lo = $wsCuentas.ListObjects.Item(1)
$lo.Resize($wsCuentas.Range("F4:H47"))
$wsCuentas.Range("H4").Value = "Tipo"
$wsCuentas.Range("H5:H47").Value = "balance"

# --- 2. Move the whole table (now F4:H47) to A1:C44 ---
$wsCuentas.Range("F4:H47").Cut($wsCuentas.Range("A1"))
$lo.Resize($wsCuentas.Range("A1:C44"))

# --- 3. Re-create the search/filter helper block at its new location E1:G6 ---
$wsCuentas.Range("E1").Value = "Search:"
$wsCuentas.Range("F1").Value = "tarjeta"
$wsCuentas.Range("E3").Value = "Cuenta"
$wsCuentas.Range("F3").Value = "Descripción"
$wsCuentas.Range("G3").Value = "Tipo"
$wsCuentas.Range("E4:G6").FormulaArray = "=FILTER(Cuentas[],(ISNUMBER(SEARCH(F1,Cuentas[Descripción]))))"

# --- 4. Swap the column-width formatting to match the new layout ---
$wsCuentas.Columns.Item(3).ColumnWidth = 8.43
$wsCuentas.Columns.Item(7).ColumnWidth = 8.43
$wsCuentas.Columns.Item(2).ColumnWidth = 39.21875
$wsCuentas.Columns.Item(6).ColumnWidth = 24.6640625

# --- 5. Rename the worksheet tab "Cuentas" -> "cuentas" ---
$wsCuentas.Name = "cuentas"

# --- 6. Update the "simple" sheet's filter helper block to include the new "Tipo" column ---
$wsSimple = $wb.Worksheets.Item("simple")
$wsSimple.Range("J6:L7").FormulaArray = "=FILTER(Cuentas[],(ISNUMBER(SEARCH(K3,Cuentas[Descripción]))))"

# --- 7. Update the "compleja" sheet's filter helper block to include the new "Tipo" column ---
$wsCompleja = $wb.Worksheets.Item("compleja")
$wsCompleja.Range("O6:Q7").FormulaArray = "=FILTER(Cuentas[],(ISNUMBER(SEARCH(P3,Cuentas[Descripción]))))"
